# Mise à jour de l'application
#
# Append the stats for the new match "N3 J11 VS ASPTT Dijon (Match arrêté 83ème)"
# (date 2025-12-13) as 14 new player rows (1000-1013) at the bottom of Feuil1,
# following the same column layout used throughout the sheet:
#   A=Match, B=Date, C=Période, E=Nom du joueur, F=Poste, G=Temps joué, H:V=stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Disable multi-threaded calculation (workbook.xml calcPr concurrentCalc="0").
try {
    $excel.MultiThreadedCalculation.Enabled = $false
} catch {
    # Not all hosts expose this switch; ignore if unavailable.
}

$firstRow = 1000
$lastRow  = 1013

# ---------------------------------------------------------------------------
# Step 1 - column G ("Temps joué" per player), written row by row, first.
# This reproduces the shared-string table order of the source file, where
# the ten new time values were entered before the (single, reused) match
# name, and a handful of rows reuse times already used by earlier matches.
# ---------------------------------------------------------------------------
$gValues = @(
    '00:08:04',
    '01:27:09',
    '01:06:18',
    '00:20:17',
    '01:27:31',
    '01:05:33',
    '01:19:39',
    '00:24:36',
    '00:18:47',
    '01:27:09',
    '00:06:46',
    '01:17:57',
    '01:26:58',
    '01:19:04'
)
for ($i = 0; $i -lt $gValues.Count; $i++) {
    $ws.Range("G" + ($firstRow + $i)).Value = $gValues[$i]
}

# ---------------------------------------------------------------------------
# Step 2 - the rest of the row data (everything except A and G), none of
# which introduces any new shared string (player names, positions and the
# "Global" period all already exist elsewhere in the sheet).
# ---------------------------------------------------------------------------
$bcdefRows = @(
    @(46004, 'Global', $null, 'Levy Ndoutoume', 'left back'),
    @(46004, 'Global', $null, 'Naim Dhib', 'center midfield'),
    @(46004, 'Global', $null, 'Mattheo Haon', 'right back'),
    @(46004, 'Global', $null, 'Karim Belmahi', 'left forward'),
    @(46004, 'Global', $null, 'Yoann Martelat', 'center midfield'),
    @(46004, 'Global', $null, 'Amir Etien', 'right forward'),
    @(46004, 'Global', $null, 'Ilan Ihaddadene', 'center midfield'),
    @(46004, 'Global', $null, 'Maé Clavel', 'left back'),
    @(46004, 'Global', $null, 'Malik Boussaid', 'right back'),
    @(46004, 'Global', $null, 'Sofiane Belle', 'left forward'),
    @(46004, 'Global', $null, 'Romain Thunet', 'center back'),
    @(46004, 'Global', $null, 'Emmanuel Valey', 'left forward'),
    @(46004, 'Global', $null, 'Naim Ighbane', 'center back'),
    @(46004, 'Global', $null, 'Yoan Zouma', 'center back')
)

$hvRows = @(
    @(0.87, 0.21, 0.65, 0.15, 0.04, 0.03, 0, 2, 6.37, 26.47, 3.42, 3, 0, 2, 0),
    @(9.43, 1.61, 7.8, 1.16, 0.41, 0.06, 0, 4, 6.44, 28.59, 4.36, 29, 2, 27, 5),
    @(8.13, 1.51, 6.6, 0.97, 0.44, 0.1, 0.02, 7, 7.19, 31.58, 4.57, 27, 4, 21, 5),
    @(2.33, 0.66, 1.67, 0.43, 0.16, 0.07, 0, 4, 6.86, 29.98, 4.46, 11, 2, 11, 2),
    @(10.97, 2.41, 8.53, 1.88, 0.5, 0.05, 0, 2, 7.44, 29.43, 4.52, 34, 3, 25, 5),
    @(6.69, 1.43, 5.24, 0.79, 0.42, 0.19, 0.04, 14, 6.04, 32.48, 4.7, 33, 4, 19, 7),
    @(10.07, 2.24, 7.81, 1.56, 0.63, 0.07, 0, 5, 7.49, 28.71, 5.02, 43, 4, 22, 4),
    @(3.34, 0.74, 2.58, 0.5, 0.19, 0.06, 0, 3, 8.13, 30.04, 4.29, 15, 3, 11, 9),
    @(2.14, 0.62, 1.51, 0.31, 0.21, 0.11, 0, 5, 6.73, 28.48, 4.28, 8, 1, 6, 4),
    @(9.32, 2.14, 7.15, 1.37, 0.58, 0.21, 0.02, 16, 6.3, 32.23, 4.19, 30, 1, 36, 9),
    @(0.7, 0.07, 0.63, 0.07, 0, 0, 0, 0, 6.16, 20.21, 3.28, 1, 0, 1, 0),
    @(9.54, 2.32, 7.2, 1.33, 0.71, 0.29, 0.01, 24, 7.26, 30.6, 4.57, 42, 8, 29, 10),
    @(9.13, 1.46, 7.66, 0.84, 0.52, 0.11, 0, 10, 6.19, 29.76, 5.11, 36, 5, 16, 5),
    @(7.75, 0.74, 7, 0.56, 0.18, 0.02, 0, 4, 5.75, 28.15, 4.59, 13, 4, 8, 3)
)

$rowCount = $bcdefRows.Count

$bcdefArr = New-Object 'object[,]' $rowCount, 5
for ($i = 0; $i -lt $rowCount; $i++) {
    for ($j = 0; $j -lt 5; $j++) {
        $bcdefArr[$i, $j] = $bcdefRows[$i][$j]
    }
}
$ws.Range("B$firstRow`:F$lastRow").Value = $bcdefArr

$hvArr = New-Object 'object[,]' $rowCount, 15
for ($i = 0; $i -lt $rowCount; $i++) {
    for ($j = 0; $j -lt 15; $j++) {
        $hvArr[$i, $j] = $hvRows[$i][$j]
    }
}
$ws.Range("H$firstRow`:V$lastRow").Value = $hvArr

# Match the date-number-format style used for column B on existing rows
# (numFmtId 14, "m/d/yyyy"), without creating a redundant duplicate style.
$ws.Range("B999").Copy() | Out-Null
$ws.Range("B$firstRow`:B$lastRow").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Step 3 - column A (match name) written last for every new row, matching
# the source file where this single new shared string was appended after
# every other new string.
# ---------------------------------------------------------------------------
$ws.Range("A$firstRow`:A$lastRow").Value = 'N3 J11 VS ASPTT Dijon (Match arrêté 83ème)'

# ---------------------------------------------------------------------------
# Reflect the updated selection/view, as recorded in the edited workbook
# (scrolled down to the newly-entered data, cell D1000 selected).
# ---------------------------------------------------------------------------
$ws.Range("D1000").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 977
$excel.ActiveWindow.ScrollColumn = 1
